# SummonCount.xlsx fix:
#  - 확률 백만분율로 수정 / 10레벨까지 확장 적용 (동료, 스킬도 적용) / 레벨별 확률 밸런싱
# Adds grade 9 & 10 rows to Equipment / Skills / Follower sheets, extends the
# tables, and repoints the JSON export paths at the new @Resources/Texts/Summon
# location.

$wb = $excel.ActiveWorkbook

$equipmentPath = 'D:\Project\TeamProject-IdleGame\IdleGame\Assets\@Resources\Texts\Summon\SummonCountEquipment.json'
$skillsPath    = 'D:\Project\TeamProject-IdleGame\IdleGame\Assets\@Resources\Texts\Summon\SummonCountSkills.json'
$followerPath  = 'D:\Project\TeamProject-IdleGame\IdleGame\Assets\@Resources\Texts\Summon\SummonCountFollower.json'

function Extend-SummonSheet {
    param($ws, $titlePath, $selectCell)

    # New grade 9 / grade 10 rows.
    $ws.Range("A11").Value = 9
    $ws.Range("B11").Value = 12000
    $ws.Range("D11").Formula = "=B11-B10"

    $ws.Range("A12").Value = 10
    $ws.Range("B12").Value = 18000
    $ws.Range("D12").Formula = "=B12-B11"

    # Grow the table/autofilter to cover the new rows.
    $lo = $ws.ListObjects.Item(1)
    $lo.Resize($ws.Range("A2:B12"))

    # Repoint the exported json path held in A1.
    $ws.Range("A1").Value = $titlePath

    $ws.Range($selectCell).Select()
}

$wsEquipment = $wb.Worksheets.Item("Equipment")
Extend-SummonSheet $wsEquipment $equipmentPath "E15"

$wsSkills = $wb.Worksheets.Item("Skills")
Extend-SummonSheet $wsSkills $skillsPath "E12"

$wsFollower = $wb.Worksheets.Item("Follower")
Extend-SummonSheet $wsFollower $followerPath "A2"

# Skills ends up the active tab (was Follower).
$wsSkills.Activate()
